$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("P2").Value = "320018556693"
$ws.Range("Q2").Value = "$18.40"
$ws.Range("R2").Value = "FAIL"

$ws.Range("P4").Value = "320018556708"
$ws.Range("Q4").Value = "$9.00"
$ws.Range("R4").Value = "FAIL"
